$wb = $excel.ActiveWorkbook

# Rename the "Hoja1" sheet to "Inventario"
$wsInventario = $wb.Worksheets.Item("Hoja1")
$wsInventario.Name = "Inventario"

# --- Cartera: drop the current selection (J15) in favor of J33; it will no
#     longer be the active/selected tab once Ventas is reactivated below. ---
$wsCartera = $wb.Worksheets.Item("Cartera")
[void]$wsCartera.Activate()
[void]$wsCartera.Range("J33").Select()

# --- Inventario (ex "Hoja1"): selection moves from A10:XFD10 to F40 ---
[void]$wsInventario.Activate()
[void]$wsInventario.Range("F40").Select()

# --- Ventas: becomes the active/selected tab again; keep the frozen-pane
#     view scrolled to the bottom selection (D1239), nudged up a couple of
#     rows so topLeftCell reflects the new scroll position. ---
$wsVentas = $wb.Worksheets.Item("Ventas")
[void]$wsVentas.Activate()
[void]$wsVentas.Range("D1239").Select()
$excel.ActiveWindow.ScrollRow = 1242
